$d = $word.ActiveDocument

# 1) Fix spacing in initials "Грачев А. В." -> "Грачев А.В."
$d.Content.Find.Execute("Грачев А. В.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Грачев А.В.", 2)

# 2) Fix spacing in initials "Буланов А. А." -> "Буланов А.А."
$d.Content.Find.Execute("Буланов А. А.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Буланов А.А.", 2)

# 3) Fix missing "м" -> "мм" in "_цапфы =  2,0977 м;" -> "_цапфы =  2,0977 мм;"
$d.Content.Find.Execute("2,0977 м;", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2,0977 мм;", 2)
